$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 162
$ws.Range("F3").Value = 94
$ws.Range("F4").Value = 395
$ws.Range("F5").Value = 937
$ws.Range("F6").Value = 4856
$ws.Range("F7").Value = 372
$ws.Range("F8").Value = 552
$ws.Range("F9").Value = 854
$ws.Range("F10").Value = 797
$ws.Range("F13").Value = 532
$ws.Range("F16").Value = 1561
$ws.Range("F17").Value = 1415
$ws.Range("F18").Value = 661
$ws.Range("F20").Value = 166
$ws.Range("F21").Value = 234
$ws.Range("F22").Value = 454
$ws.Range("F23").Value = 103
$ws.Range("F27").Value = 1351
$ws.Range("F28").Value = 126
$ws.Range("F29").Value = 70
$ws.Range("F30").Value = 8
$ws.Range("F31").Value = 179
$ws.Range("F36").Value = 244
$ws.Range("F37").Value = 530
$ws.Range("F38").Value = 67
$ws.Range("F39").Value = 3
$ws.Range("F40").Value = 9
$ws.Range("F41").Value = 49

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 125
$ws.Range("F6").Value = 93

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 162
$ws.Range("F4").Value = 94
$ws.Range("F5").Value = 395
$ws.Range("F6").Value = 937
$ws.Range("F8").Value = 4856
$ws.Range("F9").Value = 372
$ws.Range("F10").Value = 552
$ws.Range("F12").Value = 125
$ws.Range("F13").Value = 854
$ws.Range("F14").Value = 797
$ws.Range("F16").Value = 93
$ws.Range("F19").Value = 532
$ws.Range("F23").Value = 1561
$ws.Range("F24").Value = 1415
$ws.Range("F25").Value = 662
$ws.Range("F27").Value = 166
$ws.Range("F28").Value = 234
$ws.Range("F30").Value = 454
$ws.Range("F31").Value = 103
$ws.Range("F34").Value = 1351
$ws.Range("F35").Value = 126
$ws.Range("F36").Value = 70
$ws.Range("F37").Value = 8
$ws.Range("F38").Value = 179
$ws.Range("F42").Value = 244
$ws.Range("F43").Value = 530
$ws.Range("F44").Value = 67
$ws.Range("F45").Value = 3
$ws.Range("F46").Value = 9
$ws.Range("F47").Value = 49
